$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-09-09 12:10:26"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-09-09 12:10:09"
$zhcn.Range("K3").Value = "2016-09-09 12:11:08"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-09-09 12:10:26"
$dede.Range("K3").Value = "2016-09-09 12:11:33"
